# The "species" sheet had a leading column A that held a running index /
# "model id" number for each row (header text "model id" in A2, values
# 0,1,2,... below it). That whole column is removed, shifting every other
# column (module/ID/name/Yinit/Ymax/tau/type/gene name/notes) one slot to
# the left. Deleting the entire column (rather than just clearing it) is
# what drops the "model id" entry from the shared-string table and shrinks
# the sheet's used range from A1:DL37 to A1:DK37.
$wb  = $excel.ActiveWorkbook
$species   = $wb.Worksheets.Item("species")
$reactions = $wb.Worksheets.Item("reactions")

$species.Columns.Item(1).Delete()

# Reflect the new cursor / active-tab state recorded in the saved file:
# the "reactions" sheet's selection moved from C18 to H17 and it is no
# longer the active tab, while "species" becomes the active tab with its
# selection on B22 (where the "model id" column used to be is now gone,
# so B22 lines up with the old C22 cell).
$reactions.Range("H17").Select()

$species.Activate()
$species.Range("B22").Select()
